# Update column G ("K" = strikeouts) on the active worksheet to reflect
# the regenerated save_data (K instead of Strike#, std/mean, s_vals regen).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    4  = 0
    5  = 0
    6  = 1
    7  = 2
    8  = 1
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    14 = 3
    15 = 1
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 1
    21 = 0
    23 = 2
    24 = 1
    25 = 1
    26 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
